$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 345.91666
$ws.Range("I2").Value = 305.1
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 305.1
$ws.Range("L2").Value = 550
$ws.Range("M2").Value = -192.1
$ws.Range("N2").Value = -776
$ws.Range("H6").Value = 229.66667
$ws.Range("I6").Value = 255.6
$ws.Range("K6").Value = 766.8
$ws.Range("M6").Value = -654.8
$ws.Range("H9").Value = 2678.5334
$ws.Range("I9").Value = 342.5
$ws.Range("J9").Value = 4235.8887
$ws.Range("K9").Value = 342.5
$ws.Range("L9").Value = 4235.8887
$ws.Range("M9").Value = -173.5
$ws.Range("N9").Value = -4573.8887
$ws.Range("H11").Value = 668.1818
$ws.Range("I11").Value = 668.1818
$ws.Range("K11").Value = 668.1818
$ws.Range("M11").Value = -528.1818
$ws.Range("H31").Value = 63.25
$ws.Range("I31").Value = 63.25
$ws.Range("K31").Value = 189.75
$ws.Range("M31").Value = 40.25
$ws.Range("H32").Value = 1996.6666
$ws.Range("I32").Value = 1945
$ws.Range("J32").Value = 2100
$ws.Range("K32").Value = 1945
$ws.Range("L32").Value = 2100
$ws.Range("M32").Value = -1619
$ws.Range("N32").Value = -2752
$ws.Range("H33").Value = 764.4
$ws.Range("I33").Value = 192.3
$ws.Range("K33").Value = 192.3
$ws.Range("M33").Value = 36.69999999999999
$ws.Range("H40").Value = 5628.2856
$ws.Range("I40").Value = 3200
$ws.Range("J40").Value = 7449.5
$ws.Range("K40").Value = 3200
$ws.Range("L40").Value = 7449.5
$ws.Range("M40").Value = -3025
$ws.Range("N40").Value = -7799.5
$ws.Range("H55").Value = 141.5
$ws.Range("I55").Value = 83.875
$ws.Range("J55").Value = 372
$ws.Range("K55").Value = 83.875
$ws.Range("L55").Value = 372
$ws.Range("M55").Value = 130.125
$ws.Range("N55").Value = -800
$ws.Range("H64").Value = 10499.583
$ws.Range("I64").Value = 8666.223
$ws.Range("J64").Value = 15999.667
$ws.Range("K64").Value = 8666.223
$ws.Range("L64").Value = 15999.667
$ws.Range("M64").Value = -8418.223
$ws.Range("N64").Value = -16495.667
$ws.Range("H67").Value = 10499.583
$ws.Range("I67").Value = 8666.223
$ws.Range("J67").Value = 15999.667
$ws.Range("K67").Value = 8666.223
$ws.Range("L67").Value = 15999.667
$ws.Range("M67").Value = -7808.223
$ws.Range("N67").Value = -17715.667
$ws.Range("H75").Value = 51249.75
$ws.Range("J75").Value = 51249.75
$ws.Range("L75").Value = 51249.75
$ws.Range("N75").Value = -53121.75
$ws.Range("H78").Value = 51249.75
$ws.Range("J78").Value = 51249.75
$ws.Range("L78").Value = 153749.25
$ws.Range("N78").Value = -163109.25
$ws.Range("H80").Value = 3600.25
$ws.Range("I80").Value = 2133.3333
$ws.Range("J80").Value = 5067.1665
$ws.Range("K80").Value = 6399.999899999999
$ws.Range("L80").Value = 15201.4995
$ws.Range("M80").Value = -5401.999899999999
$ws.Range("N80").Value = -17197.4995
$ws.Range("H83").Value = 3600.25
$ws.Range("I83").Value = 2133.3333
$ws.Range("J83").Value = 5067.1665
$ws.Range("K83").Value = 19199.9997
$ws.Range("L83").Value = 45604.4985
$ws.Range("M83").Value = -14207.9997
$ws.Range("N83").Value = -55588.4985
$ws.Range("H88").Value = 4211.857
$ws.Range("J88").Value = 4211.857
$ws.Range("L88").Value = 4211.857
$ws.Range("N88").Value = -5023.857
$ws.Range("H91").Value = 4211.857
$ws.Range("J91").Value = 4211.857
$ws.Range("L91").Value = 4211.857
$ws.Range("N91").Value = -7019.857
$ws.Range("H99").Value = 1137
$ws.Range("I99").Value = 1265
$ws.Range("J99").Value = 881
$ws.Range("K99").Value = 3795
$ws.Range("L99").Value = 2643
$ws.Range("M99").Value = -2297
$ws.Range("N99").Value = -5639
$ws.Range("H100").Value = 2593.0417
$ws.Range("I100").Value = 2831.25
$ws.Range("J100").Value = 1402
$ws.Range("K100").Value = 2831.25
$ws.Range("L100").Value = 1402
$ws.Range("M100").Value = -2290.25
$ws.Range("N100").Value = -2484
$ws.Range("H116").Value = 6281.364
$ws.Range("I116").Value = 5752.769
$ws.Range("K116").Value = 5752.769
$ws.Range("M116").Value = -2310.769
$ws.Range("H137").Value = 2254
$ws.Range("I137").Value = 957.4286
$ws.Range("J137").Value = 3766.6667
$ws.Range("K137").Value = 2872.2858
$ws.Range("L137").Value = 11300.0001
$ws.Range("M137").Value = -322.2857999999997
$ws.Range("N137").Value = -16400.0001
$ws.Range("H138").Value = 7083.811
$ws.Range("I138").Value = 5180.5835
$ws.Range("J138").Value = 7997.36
$ws.Range("K138").Value = 15541.7505
$ws.Range("L138").Value = 23992.08
$ws.Range("M138").Value = -10401.7505
$ws.Range("N138").Value = -34272.08
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H21").Value = 141336.4
$ws.Range("I21").Value = 176045.75
$ws.Range("J21").Value = 2499
$ws.Range("K21").Value = 176045.75
$ws.Range("L21").Value = 2499
$ws.Range("M21").Value = -175671.75
$ws.Range("N21").Value = -3247
$ws.Range("H32").Value = 4623.021
$ws.Range("I32").Value = 3519.6956
$ws.Range("K32").Value = 3519.6956
$ws.Range("M32").Value = -3232.6956
$ws.Range("H61").Value = 2909.7693
$ws.Range("I61").Value = 2680.2222
$ws.Range("K61").Value = 2680.2222
$ws.Range("M61").Value = -2468.2222
$ws.Range("H63").Value = 1531.1875
$ws.Range("I63").Value = 1531.1875
$ws.Range("K63").Value = 1531.1875
$ws.Range("M63").Value = -845.1875
$ws.Range("H66").Value = 1531.1875
$ws.Range("I66").Value = 1531.1875
$ws.Range("K66").Value = 7655.9375
$ws.Range("M66").Value = -4223.9375
$ws.Range("H74").Value = 1714
$ws.Range("I74").Value = 1599.6
$ws.Range("K74").Value = 1599.6
$ws.Range("M74").Value = -725.5999999999999
$ws.Range("H77").Value = 1714
$ws.Range("I77").Value = 1599.6
$ws.Range("K77").Value = 7998
$ws.Range("M77").Value = -3630
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H82").Value = 36586.25
$ws.Range("I82").Value = 23082
$ws.Range("J82").Value = 50090.5
$ws.Range("K82").Value = 23082
$ws.Range("L82").Value = 50090.5
$ws.Range("M82").Value = -22721
$ws.Range("N82").Value = -50812.5
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H85").Value = 36586.25
$ws.Range("I85").Value = 23082
$ws.Range("J85").Value = 50090.5
$ws.Range("K85").Value = 23082
$ws.Range("L85").Value = 50090.5
$ws.Range("M85").Value = -21834
$ws.Range("N85").Value = -52586.5
$ws.Range("H97").Value = 373.85715
$ws.Range("I97").Value = 328.08334
$ws.Range("K97").Value = 328.08334
$ws.Range("M97").Value = 167.91666
$ws.Range("H102").Value = 2185.077
$ws.Range("I102").Value = 1140.7
$ws.Range("J102").Value = 5666.3335
$ws.Range("K102").Value = 1140.7
$ws.Range("L102").Value = 5666.3335
$ws.Range("M102").Value = 481.3
$ws.Range("N102").Value = -8910.333500000001
$ws.Range("H136").Value = 2909.7693
$ws.Range("I136").Value = 2680.2222
$ws.Range("K136").Value = 8040.6666
$ws.Range("M136").Value = -5490.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H86").Value = 3176.5
$ws.Range("I86").Value = 2761.6
$ws.Range("K86").Value = 2761.6
$ws.Range("M86").Value = -1638.6
$ws.Range("H87").Value = 47555
$ws.Range("J87").Value = 47555
$ws.Range("L87").Value = 47555
$ws.Range("N87").Value = -50051
$ws.Range("H89").Value = 3176.5
$ws.Range("I89").Value = 2761.6
$ws.Range("K89").Value = 13808
$ws.Range("M89").Value = -8192
$ws.Range("H90").Value = 47555
$ws.Range("J90").Value = 47555
$ws.Range("L90").Value = 142665
$ws.Range("N90").Value = -155145
$ws.Range("H95").Value = 19998
$ws.Range("J95").Value = 19998
$ws.Range("L95").Value = 19998
$ws.Range("N95").Value = -25490
$ws.Range("H99").Value = 3170
$ws.Range("I99").Value = 2595
$ws.Range("J99").Value = 4320
$ws.Range("K99").Value = 2595
$ws.Range("L99").Value = 4320
$ws.Range("M99").Value = -1097
$ws.Range("N99").Value = -7316
$ws.Range("H107").Value = 1000.3889
$ws.Range("I107").Value = 962.9375
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 962.9375
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = 957.0625
$ws.Range("N107").Value = -5140
$ws.Range("H134").Value = 2446.5454
$ws.Range("I134").Value = 2174.111
$ws.Range("J134").Value = 3672.5
$ws.Range("K134").Value = 6522.333
$ws.Range("L134").Value = 11017.5
$ws.Range("M134").Value = -3987.333
$ws.Range("N134").Value = -16087.5
$ws.Range("H137").Value = 35000
$ws.Range("J137").Value = 35000
$ws.Range("L137").Value = 35000
$ws.Range("N137").Value = -45200
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2652.2856
$ws.Range("I31").Value = 2441.75
$ws.Range("K31").Value = 2441.75
$ws.Range("M31").Value = -2146.75
$ws.Range("H34").Value = 2652.2856
$ws.Range("I34").Value = 2441.75
$ws.Range("K34").Value = 2441.75
$ws.Range("M34").Value = -2239.75
$ws.Range("H43").Value = 49000
$ws.Range("J43").Value = 49000
$ws.Range("L43").Value = 49000
$ws.Range("N43").Value = -49368
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H87").Value = 99999
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 99999
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 99999
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -102371
$ws.Range("H88").Value = 21973.666
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 25368.4
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 25368.4
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -26180.4
$ws.Range("H90").Value = 99999
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 99999
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 299997
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -311853
$ws.Range("H91").Value = 21973.666
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 25368.4
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 25368.4
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -28176.4
$ws.Range("H93").Value = 23749.5
$ws.Range("I93").Value = 23749.5
$ws.Range("K93").Value = 23749.5
$ws.Range("M93").Value = -21877.5
$ws.Range("H95").Value = 11300
$ws.Range("J95").Value = 11300
$ws.Range("L95").Value = 11300
$ws.Range("N95").Value = -16792
$ws.Range("H96").Value = 17749
$ws.Range("J96").Value = 17749
$ws.Range("L96").Value = 17749
$ws.Range("N96").Value = -23241
$ws.Range("H101").Value = 49000
$ws.Range("J101").Value = 49000
$ws.Range("L101").Value = 49000
$ws.Range("N101").Value = -55490
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -54868
$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 1500
$ws.Range("M122").Value = 950
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H141").Value = 567113.8
$ws.Range("J141").Value = 567113.8
$ws.Range("L141").Value = 567113.8
$ws.Range("N141").Value = -577473.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 239
$ws.Range("I2").Value = 281.66666
$ws.Range("J2").Value = 175
$ws.Range("K2").Value = 1689.99996
$ws.Range("L2").Value = 1050
$ws.Range("M2").Value = -1576.99996
$ws.Range("N2").Value = -1276
$ws.Range("H6").Value = 172.75
$ws.Range("I6").Value = 196
$ws.Range("J6").Value = 134
$ws.Range("K6").Value = 588
$ws.Range("L6").Value = 402
$ws.Range("M6").Value = -475
$ws.Range("N6").Value = -628
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16687449
$ws.Range("I70").Value = 20019960
$ws.Range("J70").Value = 24898
$ws.Range("K70").Value = 20019960
$ws.Range("L70").Value = 24898
$ws.Range("M70").Value = -20019690
$ws.Range("N70").Value = -25438
$ws.Range("H73").Value = 16687449
$ws.Range("I73").Value = 20019960
$ws.Range("J73").Value = 24898
$ws.Range("K73").Value = 20019960
$ws.Range("L73").Value = 24898
$ws.Range("M73").Value = -20019024
$ws.Range("N73").Value = -26770
$ws.Range("H80").Value = 3042.875
$ws.Range("I80").Value = 2760
$ws.Range("J80").Value = 3514.3333
$ws.Range("K80").Value = 2760
$ws.Range("L80").Value = 3514.3333
$ws.Range("M80").Value = -1762
$ws.Range("N80").Value = -5510.3333
$ws.Range("H83").Value = 3042.875
$ws.Range("I83").Value = 2760
$ws.Range("J83").Value = 3514.3333
$ws.Range("K83").Value = 13800
$ws.Range("L83").Value = 17571.6665
$ws.Range("M83").Value = -8808
$ws.Range("N83").Value = -27555.6665
$ws.Range("H97").Value = 2303.889
$ws.Range("I97").Value = 1623.3334
$ws.Range("J97").Value = 3665
$ws.Range("K97").Value = 1623.3334
$ws.Range("L97").Value = 3665
$ws.Range("M97").Value = -1127.3334
$ws.Range("N97").Value = -4657
$ws.Range("H132").Value = 4687.125
$ws.Range("I132").Value = 4699.6
$ws.Range("K132").Value = 14098.8
$ws.Range("M132").Value = -11568.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3699.9524
$ws.Range("I46").Value = 3114.2856
$ws.Range("J46").Value = 3992.7856
$ws.Range("K46").Value = 3114.2856
$ws.Range("L46").Value = 3992.7856
$ws.Range("M46").Value = -2926.2856
$ws.Range("N46").Value = -4368.7856
$ws.Range("H55").Value = 936.5
$ws.Range("I55").Value = 261.5
$ws.Range("J55").Value = 2624
$ws.Range("K55").Value = 261.5
$ws.Range("L55").Value = 2624
$ws.Range("M55").Value = -88.5
$ws.Range("N55").Value = -2970
$ws.Range("H74").Value = 47499.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 47499.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 47499.5
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -49495.5
$ws.Range("H77").Value = 47499.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 47499.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 142498.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -152482.5
$ws.Range("H92").Value = 24000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 24000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 24000
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -28992
$ws.Range("H93").Value = 900.3
$ws.Range("I93").Value = 900.375
$ws.Range("K93").Value = 900.375
$ws.Range("M93").Value = 347.625
$ws.Range("H118").Value = 83603
$ws.Range("J118").Value = 83603
$ws.Range("L118").Value = 83603
$ws.Range("N118").Value = -86917
$ws.Range("H136").Value = 5283.4287
$ws.Range("I136").Value = 4998.8
$ws.Range("J136").Value = 5995
$ws.Range("K136").Value = 14996.4
$ws.Range("L136").Value = 17985
$ws.Range("M136").Value = -12446.4
$ws.Range("N136").Value = -23085
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3650
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3650
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -31240
$ws.Range("H126").Value = 2866.3333
$ws.Range("J126").Value = 2799.75
$ws.Range("L126").Value = 8399.25
$ws.Range("N126").Value = -13339.25
$ws.Range("H132").Value = 2347.4783
$ws.Range("I132").Value = 2396.6316
$ws.Range("J132").Value = 2114
$ws.Range("K132").Value = 7189.8948
$ws.Range("L132").Value = 6342
$ws.Range("M132").Value = -4659.8948
$ws.Range("N132").Value = -11402
